$wb = $excel.ActiveWorkbook

# zh-cn sheet
$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($row in 4..7) {
    $wsZh.Range("E$row").Value = "ht"
    $wsZh.Range("H$row").Value = "2016-08-26 22:30:39"
}

# de-de sheet
$wsDe = $wb.Worksheets.Item("de-de")
foreach ($row in 4..7) {
    $wsDe.Range("E$row").Value = "ht"
    $wsDe.Range("H$row").Value = "2016-08-26 22:30:45"
}

# Overview sheet mirrors the de-de "Latest Handoff Datetime" in its
# "Latest HO Xliff Generate Date" column (G)
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($row in 4..7) {
    $wsOverview.Range("G$row").Value = "2016-08-26 22:30:45"
}
